$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the scraped crypto price/volume snapshot.
# Column D holds prices as literal text (e.g. "1.025", "27.324.56"),
# so force a Text number format before writing, then clear the applied
# format again so the cell keeps its original (default) style.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.324.56'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +3.93%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.836.52'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +4.10%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.025'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +3.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '319.53'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.71%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.021'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.61%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4336'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +1.76%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3713'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +2.23%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07317'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +2.23%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8743'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +3.30%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '2.038.61'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +15.86%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.31'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +5.06%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.467'
$ws.Range('D13').ClearFormats()
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.664'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +4.00%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.07146'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +4.18%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '81.97'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +4.23%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.025'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +2.49%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008980'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +3.78%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.018'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +2.05%  '
$ws.Range('E20').Value = '  +2.84%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '27.349.44'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +3.92%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.232'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +3.08%  '
$ws.Range('E23').Value = '  +0.21%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.276.72'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +14.59%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '156.38'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +3.42%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.897'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +2.59%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.51'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +2.94%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.273'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +3.80%  '
$ws.Range('E29').Value = '  +7.73%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '115.32'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.59%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08996'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.92%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.197'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +7.09%  '
$ws.Range('E33').Value = '  +4.49%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.456'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +3.54%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.842'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +4.85%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.021'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +2.94%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.147'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +5.74%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01952'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +3.61%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05253'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +2.35%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5155'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +5.06%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.799'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +8.86%  '
$ws.Range('E42').Value = '  +3.55%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.511'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +3.65%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.437'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +5.73%  '
$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '107.88'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +3.23%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.47'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +2.33%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.023'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +3.08%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.902'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +9.68%  '
$ws.Range('B49').Value = 'Decentraland'
$ws.Range('C49').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4621'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +3.31%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.663'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +3.92%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06281'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +1.78%  '
